$d = $word.ActiveDocument

# --- Change 1: "hivemind" -> "hive" + "mind" (split into two runs) ---
$r = $d.Content
$r.Find.Execute("hivemind")
$start = $r.Start
$mid = $start + 4
$r1 = $d.Range($start, $mid)
$r1.Bold = 1
$r1.Bold = 0

# --- Change 2: Story 4 (controller) paragraph replacement ---
$r2 = $d.Content
$r2.Find.Execute("This story was not able to be completed. Time that would have been allotted for meeting this story was reallocated primarily for implementing menu systems and reworking the world building system.")
$r2.Text = "This story was successfully implemented. Standard controller compatibility implemented though Unity's local interfaces was possible within the 2 hours predicted time. Currently it has been debugged with an Xbox 360 controller."

# --- Change 3: Overall feedback paragraph restructuring ---
$r3 = $d.Content
$r3.Find.Execute("like controller implementation and multiplayer")
$r3.Text = "related to multiplayer"

$r4 = $d.Content
$r4.Find.Execute("proper judgement on")
$r4.Text = "proper judgment on"

$rfull = $d.Content
$rfull.Find.Execute("This was the final sprint for our project")
$fstart = $rfull.Start

$splits = @(439, 474, 621, 629)
foreach ($off in $splits) {
    $sp = $fstart + $off
    $rs = $d.Range($fstart, $sp)
    $rs.Bold = 1
    $rs.Bold = 0
}
